# Apply the two changes described by the diff:
#  1) Remove the red highlight from the first "MainClass" run
#     (in the "Goal ::= ( ImportDeclaration )*  MainClass ( TypeDeclaration )* <EOF>" line).
#  2) Remove the stray second "MainClass | " fragment (extra tab, the
#     highlighted "MainClass" run, and the "| " run) that duplicated the
#     rule name just before the "ClassDeclaration" alternative, right
#     before the _GoBack bookmark.

$d = $word.ActiveDocument

# --- Change 1: strip the w:highlight on the first "MainClass" occurrence ---
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("MainClass")
if ($found1) {
    $rng1.HighlightColorIndex = 0
}

# --- Change 2: delete the duplicated "<tab>MainClass | " text before the
#     bookmark, right after "TypeDeclaration ::= <tab>" ---
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("TypeDeclaration ::= ")
if ($found2) {
    $afterTabs = $rng2.End
    # $afterTabs currently sits right after " ::= "; skip the first tab
    # (kept) and remove the fixed-length fragment
    # "<tab>MainClass | " (minus the already-skipped first tab) that
    # duplicates the rule name right before "ClassDeclaration" starts.
    $fragment = "`tMainClass | "
    $delStart = $afterTabs + 1
    $delEnd = $delStart + $fragment.Length
    $delRange = $d.Range($delStart, $delEnd)
    $delRange.Delete()
}
